$d = $word.ActiveDocument

# --- Change 1: merge "1 - " run with the question text run into a single run ---
$d.Content.Find.Execute("1 - ¿Que es un motor de plantillas?¿Cual es el motor de plantillas que utiliza Symfony? ", $false, $false, $false, $false, $false, $true, 1, $false, "1 - ¿Que es un motor de plantillas?¿Cual es el motor de plantillas que utiliza Symfony? ", 2)

# --- Change 2: delete empty paragraph right after "...Twig." paragraph ---
$p = $d.Paragraphs.Item(6)
$p.Range.Delete()

# --- Change 3: delete empty paragraph right before the first "Vista" (image caption) ---
$p = $d.Paragraphs.Item(14)
$p.Range.Delete()

# --- Change 4: rewrite the "Las 2 etapas..." paragraph and the two following empty
#     paragraphs with the new expanded content ---
$p7 = $d.Paragraphs.Item(7)
$rClear = $d.Range($p7.Range.Start, $p7.Range.End - 1)
$rClear.Text = ""

$rr = $d.Range($p7.Range.Start, $p7.Range.Start)
$rr.InsertAfter("Las 2 etapas se llaman “Twig for Template Designers” y “Twig for Template Developers”; donde  en “Twig for Template Designers” se describe toda la sintaxis y semántica del motor de plantillas, las cuales son simples archivos, donde se puede generar cualquier formato de texto (HMTL, CSV, XML, etc) sin la necesidad de especificar la extensión.")
$pos = $rr.End

$rr = $d.Range($pos, $pos)
$rr.InsertBreak(6)
$pos = $pos + 1

$rr = $d.Range($pos, $pos)
$rr.InsertAfter("Una plantilla contiene variables o expresiones, que se reemplazan con valores cuando se evalúa la plantilla, y etiquetas, que controlan la lógica de la plantilla.")

$p8 = $d.Paragraphs.Item(8)
$rClear = $d.Range($p8.Range.Start, $p8.Range.End - 1)
$rClear.Text = "Los designer pueden usar tipos de delimitadores: {% ...%} y {{...}}. El primero se utiliza para ejecutar sentencias como bucles for, el segundo genera el resultado de una expresión."

$p9 = $d.Paragraphs.Item(9)
$rClear = $d.Range($p9.Range.Start, $p9.Range.End - 1)
$rClear.Text = "Twig for Template Developers esta destinado con su API, al desarrollo de interfaz de plantillas en la aplicación."

# --- Change 5: merge "Podemos usar " + "JavaScript" + " de 2 maneras..." into a single run ---
$d.Content.Find.Execute("Podemos usar JavaScript de 2 maneras en Symfony. Una usando la biblioteca que incluye, llamada ", $false, $false, $false, $false, $false, $true, 1, $false, "Podemos usar JavaScript de 2 maneras en Symfony. Una usando la biblioteca que incluye, llamada ", 2)

# --- Change 6: merge "Sino también se puede trabajar con " + "JavaScript" + " de la manera..." into a single run ---
$d.Content.Find.Execute("Sino también se puede trabajar con JavaScript de la manera tradicional bajando desde su web los archivos o bien haciendo uso a un servidor CDN externo.", $false, $false, $false, $false, $false, $true, 1, $false, "Sino también se puede trabajar con JavaScript de la manera tradicional bajando desde su web los archivos o bien haciendo uso a un servidor CDN externo.", 2)
